# Refresh the "cryptos" price table (GitHub Actions scheduled update).
# Updates Price/Volume(1h) text for every row, plus a few rows whose
# Coin/Link also shifted (rank churn: Polkadot<->WrappedEther,
# MXToken<->TrustWalletToken swapped position; BabyDogeCoin dropped off
# and Algorand/Cronos/USDD shifted up with EnergySwap newly appearing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a blank, default-styled cell as a format donor so that writing
# numeric-looking text (e.g. "1.00", "0.0510") into column D does not
# leave the cell tagged with a stray "@" text-format style - we want it
# to match the plain, unstyled inline-string cells already in the sheet.
$fmtDonor = $ws.Range("H1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.713.53"
$fmtDonor.Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.29"
$fmtDonor.Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.44"
$fmtDonor.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$fmtDonor.Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$fmtDonor.Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.57"
$fmtDonor.Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.49"
$fmtDonor.Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.618.55"
$fmtDonor.Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.05"
$fmtDonor.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.687.01"
$fmtDonor.Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.24"
$fmtDonor.Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("E19").Value = "  +4.27%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "209.34"
$fmtDonor.Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$fmtDonor.Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.65"
$fmtDonor.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.12"
$fmtDonor.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.294.77"
$fmtDonor.Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("E35").Value = "  -5.30%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +20.62%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.786"
$fmtDonor.Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$fmtDonor.Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.33"
$fmtDonor.Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.737.16"
$fmtDonor.Copy()
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.20"
$fmtDonor.Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.101"
$fmtDonor.Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$fmtDonor.Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$fmtDonor.Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.39"
$fmtDonor.Copy()
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").Value = "  -0.99%  "

$excel.CutCopyMode = $false
